# "new sites to yandex"
# Populate Лист1!A2:F4 with three new site-check rows.
#
# Columns:
#   A - timestamp string ("MM.DD.YYYY HH:MM:SS")
#   B - signed delta, stored as text (uses a real U+2212 MINUS SIGN, not a hyphen)
#   C - numeric 0
#   D/E/F - numeric-looking codes that must stay TEXT ("19", "24", "25")
#
# Columns B, D, E and F hold digit-only-looking strings, and plain
# `.Value = "19"` auto-coerces those to numbers. Toggling the cell's
# NumberFormat to Text ("@") before the write — then back to "General" —
# forces the engine to keep the literal string without leaving a visible
# number-format behind on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.NumberFormat = "General"
}

$rows = @(
    @{ Row = 2; Time = "11.25.2019 14:35:18" },
    @{ Row = 3; Time = "11.25.2019 14:36:44" },
    @{ Row = 4; Time = "11.25.2019 14:37:01" }
)

foreach ($r in $rows) {
    $row = $r.Row

    # A: timestamp (plain text — dotted date format isn't recognised as a date)
    $ws.Range("A$row").Value = $r.Time

    # B: "−6" — U+2212 MINUS SIGN followed by 6, kept as text
    Set-TextValue $ws.Range("B$row") "−6"

    # C: numeric 0
    $ws.Range("C$row").Value = 0

    # D, E, F: numeric-looking codes stored as text
    Set-TextValue $ws.Range("D$row") "19"
    Set-TextValue $ws.Range("E$row") "24"
    Set-TextValue $ws.Range("F$row") "25"
}
